$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row's label cell (A4) into A5
# so the new row's year label cell gets the same style as the others.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 30.236
$ws.Range("C5").Value = 40.313
$ws.Range("D5").Value = 47.239
$ws.Range("E5").Value = 36.415
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 34.579
$ws.Range("H5").Value = 31.783
$ws.Range("I5").Value = 33.907
$ws.Range("J5").Value = 29.319
$ws.Range("K5").Value = 47.77
$ws.Range("L5").Value = 30.296
$ws.Range("M5").Value = 34.312
$ws.Range("N5").Value = 17.301
$ws.Range("O5").Value = 7.547
$ws.Range("P5").Value = 33.537
$ws.Range("Q5").Value = 24.635
$ws.Range("R5").Value = 7.161
$ws.Range("S5").Value = 18.674
$ws.Range("T5").Value = 31.393
$ws.Range("U5").Value = 5.983
$ws.Range("V5").Value = 29.863
$ws.Range("W5").Value = 42.063
$ws.Range("X5").Value = 7.34
$ws.Range("Y5").Value = 13.353
$ws.Range("Z5").Value = 7.298
$ws.Range("AA5").Value = 6.017
$ws.Range("AB5").Value = 40.035
$ws.Range("AC5").Value = 21.562
$ws.Range("AD5").Value = 25.034
$ws.Range("AE5").Value = 10.49
$ws.Range("AF5").Value = 23.705
$ws.Range("AG5").Value = 19.759
$ws.Range("AH5").Value = 41.262
$ws.Range("AI5").Value = 36.303
$ws.Range("AJ5").Value = 28.152
$ws.Range("AK5").Value = 44.599
$ws.Range("AL5").Value = 8.986000000000001
$ws.Range("AM5").Value = 15.204
$ws.Range("AN5").Value = 28.018
$ws.Range("AO5").Value = 31.456
$ws.Range("AP5").Value = 21.195
$ws.Range("AQ5").Value = 11.815
$ws.Range("AR5").Value = 46.15
$ws.Range("AS5").Value = 19.915
$ws.Range("AT5").Value = 8.548
